# Consumo.xlsx edit: update the "Produccion - Expedición" building section.
# - C40 header text changes from "Expedición" to "Produccion - Expedición"
# - C44 switches from "Acces point" to "Switch"
# - A new "Acces point" line (row 45) is added with Cantidad=1, Consumo/U=35
# - The old "0,85*" transformer-sizing block (Pot Transf / Transf KVA / Factor
#   utliz) is removed; the Total:/C-Resguardo/A rows shift down by one row to
#   make room for the new "Acces point" line.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$xlFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# --- Section header: "Expedición" -> "Produccion - Expedición" -------------
$ws.Range("C40").Value = "Produccion - Expedición"

# --- Row 44: "Acces point" -> "Switch" --------------------------------------
$ws.Range("C44").Value = "Switch"

# --- New row 45: "Acces point" line (was blank, only M45 had the note) -----
# Copy the Cantidad/Consumo-U/Subt. formatting down from row 44 first.
$ws.Range("C44").Copy()
$ws.Range("C45").PasteSpecial($xlFormats)
$ws.Range("D44").Copy()
$ws.Range("D45").PasteSpecial($xlFormats)
$ws.Range("E44").Copy()
$ws.Range("E45").PasteSpecial($xlFormats)
$ws.Range("F44").Copy()
$ws.Range("F45").PasteSpecial($xlFormats)

$ws.Range("C45").Value = "Acces point"
$ws.Range("D45").Value = 1
$ws.Range("E45").Value = 35
$ws.Range("F45").Formula = "=D45*E45"

# --- Row 46 becomes the new blank spacer row (old "Total:" row content is --
# --- dropped/shifted) - keep the cells present but empty, matching the ----
# --- formatting that the (now gone) blank row 45 used to have. ------------
$ws.Range("C44").Copy()
$ws.Range("C46").PasteSpecial($xlFormats)
$ws.Range("D44").Copy()
$ws.Range("D46").PasteSpecial($xlFormats)
$ws.Range("E44").Copy()
$ws.Range("E46").PasteSpecial($xlFormats)
$ws.Range("F44").Copy()
$ws.Range("F46").PasteSpecial($xlFormats)
$ws.Range("C46").ClearContents()
$ws.Range("D46").ClearContents()
$ws.Range("E46").ClearContents()
$ws.Range("F46").ClearContents()
$ws.Range("K46").Clear()

# --- Row 47: Total: = SUM(F42:F46) (same look as the old row-46 totals) ---
$ws.Range("E22").Copy()
$ws.Range("E47").PasteSpecial($xlFormats)
$ws.Range("F22").Copy()
$ws.Range("F47").PasteSpecial($xlFormats)
$ws.Range("E47").Value = "Total:"
$ws.Range("F47").Formula = "=SUM(F42:F46)"

# --- Drop the old "C/Resguardo" row 47 extras (Pot Transf/0,85*/75) --------
$ws.Range("H47").Clear()
$ws.Range("I47").Clear()
$ws.Range("K47").Clear()
$ws.Range("M47").Clear()

# --- Row 48: C/Resguardo = F47 + F47*I41 (same look as old row-47) --------
$ws.Range("E23").Copy()
$ws.Range("E48").PasteSpecial($xlFormats)
$ws.Range("F23").Copy()
$ws.Range("F48").PasteSpecial($xlFormats)
$ws.Range("E48").Value = "C/Resguardo"
$ws.Range("F48").Formula = "=F47+(F47*`$I`$41)"

# --- Drop the old "Factor utliz." row 48 extras -----------------------------
$ws.Range("H48").Clear()
$ws.Range("I48").Clear()

# --- Row 49: A = F48 / L41 (same look as old row-48) ------------------------
$ws.Range("E24").Copy()
$ws.Range("E49").PasteSpecial($xlFormats)
$ws.Range("F24").Copy()
$ws.Range("F49").PasteSpecial($xlFormats)
$ws.Range("E49").Value = "A"
$ws.Range("F49").Formula = "=F48/`$L`$41"

# --- I42 total consumption formula now points at the shifted F48 (was F47) -
$ws.Range("I42").Formula = "=F13+L10+F23+L23+F35+L35+F48+F57+L57"
$ws.Range("I42").Style = "Normal"

# --- Sheet view: scroll position / selection match the saved state ---------
$ws.Range("C38").Select()
$excel.ActiveWindow.ScrollRow = 25
